$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '42.820.85'
Set-TextValue $ws 'E2' '  -1.15%  '
Set-TextValue $ws 'D3' '2.238.76'
Set-TextValue $ws 'E3' '  -1.62%  '
Set-TextValue $ws 'E4' '  +0.34%  '
Set-TextValue $ws 'D5' '115.28'
Set-TextValue $ws 'D6' '284.23'
Set-TextValue $ws 'E6' '  +7.67%  '
Set-TextValue $ws 'D7' '0.627'
Set-TextValue $ws 'E7' '  -2.40%  '
Set-TextValue $ws 'E8' '  +0.22%  '
Set-TextValue $ws 'D9' '0.613'
Set-TextValue $ws 'E9' '  +0.96%  '
Set-TextValue $ws 'D10' '46.61'
Set-TextValue $ws 'E10' '  +0.10%  '
Set-TextValue $ws 'E11' '  -0.60%  '
Set-TextValue $ws 'E12' '  -0.41%  '
Set-TextValue $ws 'E13' '  -2.74%  '
Set-TextValue $ws 'D14' '15.36'
Set-TextValue $ws 'E14' '  +0.27%  '
Set-TextValue $ws 'D15' '0.884'
Set-TextValue $ws 'E15' '  +2.94%  '
Set-TextValue $ws 'D16' '2.576.07'
Set-TextValue $ws 'E16' '  -1.64%  '
Set-TextValue $ws 'D17' '2.249.41'
Set-TextValue $ws 'E17' '  -1.22%  '
Set-TextValue $ws 'D18' '42.815.12'
Set-TextValue $ws 'E18' '  -0.71%  '
Set-TextValue $ws 'D19' '0.0000107'
Set-TextValue $ws 'E19' '  -0.70%  '
Set-TextValue $ws 'D20' '6.87'
Set-TextValue $ws 'E20' '  +2.14%  '
Set-TextValue $ws 'D21' '72.30'
Set-TextValue $ws 'E21' '  +0.22%  '
Set-TextValue $ws 'B22' 'ImmutableX'
Set-TextValue $ws 'C22' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws 'D22' '2.36'
Set-TextValue $ws 'E22' '  -2.82%  '
Set-TextValue $ws 'B23' 'PancakeSwap'
Set-TextValue $ws 'C23' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws 'D23' '3.16'
Set-TextValue $ws 'E23' '  +10.46%  '
Set-TextValue $ws 'D24' '231.90'
Set-TextValue $ws 'E24' '  -0.93%  '
Set-TextValue $ws 'D25' '9.27'
Set-TextValue $ws 'E25' '  -1.11%  '
Set-TextValue $ws 'D26' '12.05'
Set-TextValue $ws 'E26' '  +6.53%  '
Set-TextValue $ws 'E27' '  -1.67%  '
Set-TextValue $ws 'B28' 'InjectiveProtocol'
Set-TextValue $ws 'C28' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws 'D28' '40.41'
Set-TextValue $ws 'E28' '  -1.66%  '
Set-TextValue $ws 'B29' 'Toncoin'
Set-TextValue $ws 'C29' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws 'D29' '2.24'
Set-TextValue $ws 'E29' '  -0.28%  '
Set-TextValue $ws 'B30' 'WEMIXToken'
Set-TextValue $ws 'C30' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws 'D30' '3.29'
Set-TextValue $ws 'E30' '  -1.72%  '
Set-TextValue $ws 'B31' 'Monero'
Set-TextValue $ws 'C31' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws 'D31' '175.07'
Set-TextValue $ws 'E31' '  +0.91%  '
Set-TextValue $ws 'B32' 'EthereumClassic'
Set-TextValue $ws 'C32' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws 'D32' '21.17'
Set-TextValue $ws 'E32' '  -1.31%  '
Set-TextValue $ws 'B33' 'Hedera'
Set-TextValue $ws 'C33' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws 'D33' '0.0903'
Set-TextValue $ws 'E33' '  +0.95%  '
Set-TextValue $ws 'B34' 'NEARProtocol'
Set-TextValue $ws 'C34' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws 'D34' '4.61'
Set-TextValue $ws 'E34' '  +18.98%  '
Set-TextValue $ws 'B35' 'Filecoin'
Set-TextValue $ws 'C35' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws 'D35' '5.59'
Set-TextValue $ws 'E35' '  -1.01%  '
Set-TextValue $ws 'B36' 'Stellar'
Set-TextValue $ws 'C36' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws 'D36' '0.128'
Set-TextValue $ws 'E36' '  -2.47%  '
Set-TextValue $ws 'B37' 'RenderToken'
Set-TextValue $ws 'C37' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D37' '4.66'
Set-TextValue $ws 'E37' '  -0.33%  '
Set-TextValue $ws 'B38' 'VeChain'
Set-TextValue $ws 'C38' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D38' '0.0372'
Set-TextValue $ws 'E38' '  -2.13%  '
Set-TextValue $ws 'B39' 'Kaspa'
Set-TextValue $ws 'C39' 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws 'D39' '0.106'
Set-TextValue $ws 'E39' '  +1.91%  '
Set-TextValue $ws 'B40' 'LidoDAOToken'
Set-TextValue $ws 'C40' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws 'D40' '2.61'
Set-TextValue $ws 'E40' '  +1.41%  '
Set-TextValue $ws 'B41' 'MultiversX'
Set-TextValue $ws 'C41' 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextValue $ws 'D41' '72.47'
Set-TextValue $ws 'E41' '  -3.09%  '
Set-TextValue $ws 'B42' 'Celestia'
Set-TextValue $ws 'C42' 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws 'D42' '13.49'
Set-TextValue $ws 'E42' '  -5.11%  '
Set-TextValue $ws 'B43' 'Algorand'
Set-TextValue $ws 'C43' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws 'D43' '0.234'
Set-TextValue $ws 'E43' '  -0.46%  '
Set-TextValue $ws 'B44' 'FirstDigitalUSD'
Set-TextValue $ws 'C44' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws 'D44' '1.00'
Set-TextValue $ws 'E44' '  +0.40%  '
Set-TextValue $ws 'B45' 'ARBITRUM'
Set-TextValue $ws 'C45' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws 'D45' '1.34'
Set-TextValue $ws 'E45' '  -1.40%  '
Set-TextValue $ws 'B46' 'THORChain'
Set-TextValue $ws 'C46' 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws 'D46' '5.58'
Set-TextValue $ws 'E46' '  -8.12%  '
Set-TextValue $ws 'B47' 'TrustWalletToken'
Set-TextValue $ws 'C47' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws 'D47' '1.29'
Set-TextValue $ws 'E47' '  +1.31%  '
Set-TextValue $ws 'B48' 'FraxShare'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws 'D48' '8.54'
Set-TextValue $ws 'E48' '  -0.04%  '
Set-TextValue $ws 'B49' 'TheSandbox'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws 'D49' '0.649'
Set-TextValue $ws 'E49' '  +8.84%  '
Set-TextValue $ws 'B50' 'Cronos'
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws 'D50' '0.0989'
Set-TextValue $ws 'E50' '  -0.15%  '
Set-TextValue $ws 'B51' 'Aave'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D51' '100.97'
Set-TextValue $ws 'E51' '  +0.66%  '
